# Insert a new weekly price record as row 46, pushing the existing
# rows 46-53 down to 47-54 (the workbook keeps records newest-first).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 46:53 down by one row to make room for the new record.
$ws.Rows("46").Insert()

# Populate the newly inserted row 46 with the new record's data.
$ws.Range("A46").Value = 10
$ws.Range("B46").Value = "Vega Modelo de Temuco"
$ws.Range("C46").Value = "La Araucanía"
$ws.Range("D46").Value = 44505
$ws.Range("D46").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E46").Value = 9
$ws.Range("F46").Value = 100112022
$ws.Range("G46").Value = "Arveja Verde"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 125
$ws.Range("K46").Value = 16000
$ws.Range("L46").Value = 16000
$ws.Range("M46").Value = 16000
$ws.Range("N46").Value = "$/saco 25 kilos"
$ws.Range("O46").Value = "Provincia de Limarí"
$ws.Range("P46").Value = 640
$ws.Range("Q46").Value = 25
$ws.Range("R46").Value = "Hortaliza"
